{"js": "// The document contains several repeated \"product listing\" blocks, each\n// one being a run of paragraphs describing a single marketplace offer\n// (Modelo / URL / Nome / Pre\u00e7o / Pre\u00e7o Previsto / Loja / Tipo / Lugar /\n// Cupom / a dashed separator) followed by one blank paragraph.\n//\n// This edit drops the listing blocks for the stores \"DIGITALSHOP SC\",\n// \"BASSAN\", \"Best Online\" and \"RENOV VENDAS_ONLINE\", keeping the\n// \"ULTRAFER FERRAMENTAS\" and \"Radical Som\" blocks untouched.\n\nconst paras = context.document.body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst items = paras.items;\nconst total = items.length;\n\n// Locate block boundaries: each block begins at a \"Modelo:\" paragraph and\n// runs up to (but not including) the next \"Modelo:\" paragraph, or the end\n// of the body for the last block.\nconst starts = [];\nfor (let i = 0; i < total; i++) {\n  if (items[i].text.indexOf(\"Modelo:\") === 0) {\n    starts.push(i);\n  }\n}\n\nconst storesToRemove = [\n  \"Loja: DIGITALSHOP SC\",\n  \"Loja: BASSAN\",\n  \"Loja: Best Online\",\n  \"Loja: RENOV VENDAS_ONLINE\"\n];\n\nfunction blockHasTargetStore(blockStart, blockEnd) {\n  for (let i = blockStart; i <= blockEnd; i++) {\n    const t = items[i].text;\n    for (const store of storesToRemove) {\n      if (t.indexOf(store) === 0) return true;\n    }\n  }\n  return false;\n}\n\n// Collect the [start, end] (inclusive) ranges of paragraph indices that\n// must be removed.\nconst rangesToDelete = [];\nfor (let b = 0; b < starts.length; b++) {\n  const blockStart = starts[b];\n  const blockEnd = (b + 1 < starts.length) ? starts[b + 1] - 1 : total - 1;\n  if (blockHasTargetStore(blockStart, blockEnd)) {\n    rangesToDelete.push([blockStart, blockEnd]);\n  }\n}\n\n// The very last paragraph of the body can never be truly deleted (Word\n// requires the body to end with a paragraph mark); Office.js silently\n// no-ops that particular delete(). If the last range to remove reaches\n// the final paragraph, trim the range by one and instead delete the\n// (now redundant) blank paragraph that currently terminates the last\n// *kept* block, once all other removals are done.\nlet reclaimTrailingBlank = false;\nif (rangesToDelete.length > 0) {\n  const lastRange = rangesToDelete[rangesToDelete.length - 1];\n  if (lastRange[1] === total - 1) {\n    lastRange[1] = total - 2; // leave the unremovable final paragraph alone\n    reclaimTrailingBlank = true;\n  }\n}\n\n// Delete paragraph-by-paragraph, from the highest index down to the\n// lowest, so indices already visited stay valid for items not yet\n// processed.\nfor (let r = rangesToDelete.length - 1; r >= 0; r--) {\n  const [start, end] = rangesToDelete[r];\n  for (let i = end; i >= start; i--) {\n    items[i].delete();\n  }\n}\nawait context.sync();\n\nif (reclaimTrailingBlank) {\n  // Find the index of the blank paragraph directly preceding the\n  // (unremovable) final paragraph of the body and delete it, so the\n  // document ends with exactly one blank paragraph, matching the\n  // original separator pattern.\n  const paras2 = context.document.body.paragraphs;\n  paras2.load(\"items/text\");\n  await context.sync();\n  const items2 = paras2.items;\n  const n = items2.length;\n  if (n >= 2 && items2[n - 1].text === \"\" && items2[n - 2].text === \"\") {\n    items2[n - 2].delete();\n    await context.sync();\n  }\n}\n", "ps1": "# The document contains several repeated \"product listing\" blocks, each\n# one being a run of paragraphs describing a single marketplace offer\n# (Modelo / URL / Nome / Pre\u00e7o / Pre\u00e7o Previsto / Loja / Tipo / Lugar /\n# Cupom / a dashed separator) followed by one blank paragraph.\n#\n# This edit drops the listing blocks for the stores \"DIGITALSHOP SC\",\n# \"BASSAN\", \"Best Online\" and \"RENOV VENDAS_ONLINE\", keeping the\n# \"ULTRAFER FERRAMENTAS\" and \"Radical Som\" blocks untouched.\n\n$d = $word.ActiveDocument\n$n = $d.Paragraphs.Count\n\n# Collect 1-indexed positions where each block starts (\"Modelo:\" paragraphs).\n$starts = @()\nfor ($i = 1; $i -le $n; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.StartsWith(\"Modelo:\")) {\n        $starts += $i\n    }\n}\n\n$storesToRemove = @(\"Loja: DIGITALSHOP SC\", \"Loja: BASSAN\", \"Loja: Best Online\", \"Loja: RENOV VENDAS_ONLINE\")\n\n# Determine the [start, end] (1-indexed, inclusive) paragraph ranges that\n# belong to blocks which must be removed. Parallel arrays are used (rather\n# than an array-of-arrays) since element-wise re-assignment into nested\n# arrays is unreliable in this interpreter.\n$delStarts = @()\n$delEnds = @()\nfor ($b = 0; $b -lt $starts.Count; $b++) {\n    $blockStart = $starts[$b]\n    if ($b + 1 -lt $starts.Count) {\n        $blockEnd = $starts[$b + 1] - 1\n    } else {\n        $blockEnd = $n\n    }\n\n    $hasTarget = $false\n    for ($i = $blockStart; $i -le $blockEnd; $i++) {\n        $t = $d.Paragraphs.Item($i).Range.Text\n        foreach ($store in $storesToRemove) {\n            if ($t.StartsWith($store)) {\n                $hasTarget = $true\n                break\n            }\n        }\n        if ($hasTarget) { break }\n    }\n\n    if ($hasTarget) {\n        $delStarts += $blockStart\n        $delEnds += $blockEnd\n    }\n}\n\n# The very last paragraph of the body can never truly be deleted (Word\n# requires the body's main story to end with a paragraph mark); deleting\n# it is silently a no-op. If the last range to remove reaches the final\n# paragraph, trim the range by one and instead delete the (now redundant)\n# blank paragraph that currently terminates the last *kept* block, once\n# all other removals are done.\n$reclaimTrailingBlank = $false\nif ($delStarts.Count -gt 0) {\n    $lastIdx = $delStarts.Count - 1\n    if ($delEnds[$lastIdx] -eq $n) {\n        $delEnds[$lastIdx] = $n - 1\n        $reclaimTrailingBlank = $true\n    }\n}\n\n# Delete paragraph-by-paragraph, from the highest index down to the\n# lowest, so indices already visited stay valid for paragraphs not yet\n# processed.\nfor ($r = $delStarts.Count - 1; $r -ge 0; $r--) {\n    $start = $delStarts[$r]\n    $end = $delEnds[$r]\n    for ($i = $end; $i -ge $start; $i--) {\n        $d.Paragraphs.Item($i).Range.Delete()\n    }\n}\n\nif ($reclaimTrailingBlank) {\n    # Find the blank paragraph directly preceding the (unremovable) final\n    # paragraph of the document and delete it, so the document ends with\n    # exactly one blank paragraph, matching the original separator\n    # pattern.\n    $m = $d.Paragraphs.Count\n    if ($m -ge 2) {\n        $lastText = $d.Paragraphs.Item($m).Range.Text\n        $secondLastText = $d.Paragraphs.Item($m - 1).Range.Text\n        if ($lastText.Trim() -eq \"\" -and $secondLastText.Trim() -eq \"\") {\n            $d.Paragraphs.Item($m - 1).Range.Delete()\n        }\n    }\n}\n"}
